$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.957.78'
$ws.Range('E2').Value = '  -1.91%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.275.13'
$ws.Range('E3').Value = '  -1.44%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '186.26'
$ws.Range('E5').Value = '  -0.90%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '582.00'
$ws.Range('E6').Value = '  -1.80%  '

$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('E8').Value = '  -0.95%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '3.272.81'
$ws.Range('E9').Value = '  -1.55%  '

$ws.Range('E10').Value = '  -4.39%  '

$ws.Range('E11').Value = '  -2.12%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.412'
$ws.Range('E12').Value = '  -3.51%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.841.76'
$ws.Range('E13').Value = '  -1.27%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '27.62'
$ws.Range('E15').Value = '  -5.72%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '67.974.51'
$ws.Range('E16').Value = '  -1.85%  '

$ws.Range('E17').Value = '  -3.35%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.272.43'
$ws.Range('E18').Value = '  -0.73%  '

$ws.Range('E19').Value = '  -3.21%  '

$ws.Range('E20').Value = '  -1.98%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '398.55'
$ws.Range('E21').Value = '  +2.25%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.64'
$ws.Range('E22').Value = '  -2.94%  '

$ws.Range('E23').Value = '  +0.11%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '71.47'
$ws.Range('E24').Value = '  -0.70%  '

$ws.Range('E25').Value = '  -2.01%  '

$ws.Range('E26').Value = '  -5.11%  '

$ws.Range('E27').Value = '  -1.79%  '

$ws.Range('E28').Value = '  -4.44%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +0.74%  '

$ws.Range('E30').Value = '  -3.06%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '22.73'
$ws.Range('E31').Value = '  -2.32%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.52'
$ws.Range('E32').Value = '  -6.90%  '

$ws.Range('E33').Value = '  -4.48%  '

$ws.Range('E34').Value = '  -6.13%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '162.81'
$ws.Range('E36').Value = '  -0.47%  '

$ws.Range('E37').Value = '  -6.22%  '

$ws.Range('E38').Value = '  -1.42%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '26.91'
$ws.Range('E39').Value = '  -0.71%  '

$ws.Range('E40').Value = '  -4.05%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.55'
$ws.Range('E41').Value = '  -2.43%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.49'
$ws.Range('E42').Value = '  -4.09%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.678.41'
$ws.Range('E43').Value = '  +0.23%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0686'
$ws.Range('E44').Value = '  -2.13%  '

$ws.Range('E45').Value = '  -1.97%  '

$ws.Range('E46').Value = '  -9.25%  '

$ws.Range('E47').Value = '  -4.80%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '333.74'
$ws.Range('E48').Value = '  -3.15%  '

$ws.Range('E49').Value = '  -3.83%  '

$ws.Range('E50').Value = '  +0.37%  '

$ws.Range('E51').Value = '  -1.71%  '
